$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$updates = @(
    @{Row=2; C=0.0096522889319272179; D=-0.32553150662035485; E=-0.28769508267065941},
    @{Row=3; C=0.0029091036130235001; D=0.15055738047695383; E=0.16196090127831164},
    @{Row=4; C=0.011865318626738435; D=0.27510672571521716; E=0.32161810027369259},
    @{Row=5; C=0.0041453603647392447; D=-0.057504134552584928; E=-0.041254551940040123},
    @{Row=6; C=0.001957499364843861; D=0.056260221038670698; E=0.063933509931380791},
    @{Row=7; C=0.005136820117669855; D=0.039385056061103277; E=0.059521105891142359},
    @{Row=8; C=0.008544828282437994; D=-0.34094209854021967; E=-0.3074468573853199},
    @{Row=9; C=0.0054383528928735264; D=0.26978609601034459; E=0.2911041302190171},
    @{Row=10; C=0.010867165675226918; D=0.22745373700679503; E=0.27005240873318803},
    @{Row=11; C=0.0052338094410633142; D=-0.12950052093995848; E=-0.10898427834394836},
    @{Row=12; C=0.0033019628757557876; D=0.13795541615447079; E=0.15089892740853392},
    @{Row=13; C=0.0075182270282111961; D=0.068315730098452099; E=0.097786762878606964},
    @{Row=14; C=0.0064734294450451814; D=-0.29328548206860938; E=-0.26791000661210884},
    @{Row=15; C=0.0041958268841884531; D=0.34864891634260747; E=0.36509631922593438},
    @{Row=16; C=0.0098858236024648456; D=0.12642538563488373; E=0.16517725221831417},
    @{Row=17; C=0.0044314336232968718; D=-0.12136378378363048; E=-0.10399280987113667},
    @{Row=18; C=0.002884858250504638; D=0.17212531745349086; E=0.1834338017208183},
    @{Row=19; C=0.0066634238438001286; D=0.052089579464926997; E=0.07820983119340176},
    @{Row=20; C=0.0055594996467473717; D=-0.22546853378525081; E=-0.20367561118771768},
    @{Row=21; C=0.0037179180867323038; D=0.35879433986256121; E=0.37336836742555052},
    @{Row=22; C=0.009669093061444969; D=0.025646009007709795; E=0.063548304189920496},
    @{Row=23; C=0.0039825747662933951; D=-0.11858743802666655; E=-0.10297596592740271},
    @{Row=24; C=0.0038922727760216129; D=0.21834323328984739; E=0.23360072659791126},
    @{Row=25; C=0.0074026955164285379; D=0.033761308412812538; E=0.062779464077377348},
    @{Row=26; C=0.0040629196301676506; D=-0.16021805603033132; E=-0.14429164202794095},
    @{Row=27; C=0.0037406281749947073; D=0.31411993262850835; E=0.32878298244658055},
    @{Row=28; C=0.0080267148361303203; D=-0.050975161863114585; E=-0.019510895966704175},
    @{Row=29; C=0.0038117547411906849; D=-0.10700690761232302; E=-0.092065040533021822},
    @{Row=30; C=0.0042022900365883773; D=0.23948592383338543; E=0.25595866760067237},
    @{Row=31; C=0.0070928342305495969; D=0.0123372278930733; E=0.040140744510536605},
    @{Row=32; C=0.0037057912240009727; D=-0.11615061601475497; E=-0.10162412506434722},
    @{Row=33; C=0.0037582253378713151; D=0.25295691653013685; E=0.26768894622641271},
    @{Row=34; C=0.0090413200497756214; D=-0.1017059050706154; E=-0.066264444409751336},
    @{Row=35; C=0.0038207979689629376; D=-0.093665741584812515; E=-0.078688425554433547},
    @{Row=36; C=0.00452412215998939; D=0.24879776355428809; E=0.26653207138752677},
    @{Row=37; C=0.0085526401206023322; D=-0.014896661957128222; E=0.018629212747816717},
    @{Row=38; C=0.0033064515094659371; D=-0.091551123310067248; E=-0.07859002134103564},
    @{Row=39; C=0.0056157210692007612; D=0.16502732015076965; E=0.18704062752853512},
    @{Row=40; C=0.0073658551506586289; D=-0.10663922487878282; E=-0.077765491384283916},
    @{Row=41; C=0.004859627110094928; D=-0.086530474371517146; E=-0.067481005750321121},
    @{Row=42; C=0.0045116389640615017; D=0.25298272033090868; E=0.27066809472877601},
    @{Row=43; C=0.0090136393511639767; D=-0.041680508818381798; E=-0.0063475427095023677},
    @{Row=44; C=0.003850470695556339; D=-0.078562009580718545; E=-0.063468383325804956},
    @{Row=45; C=0.0040539359787561896; D=0.059723960346825679; E=0.075615158946339328},
    @{Row=46; C=0.0072851277209446466; D=-0.10248836166002975; E=-0.073931075101233756},
    @{Row=47; C=0.0039345779493968353; D=-0.06157398886788451; E=-0.046150661627613873},
    @{Row=48; C=0.0049222113608480798; D=0.20109187944772941; E=0.22038667485921079},
    @{Row=49; C=0.010301668189589282; D=-0.050790091891229155; E=-0.010408124205704247},
    @{Row=50; C=0.0051317303865732189; D=-0.088790831655165436; E=-0.068674740241897847},
    @{Row=51; C=0.0045351277903548657; D=-0.058414879780231063; E=-0.040637436631556972},
    @{Row=52; C=0.008182334415981523; D=-0.036022935216160304; E=-0.0039486494125911542},
    @{Row=53; C=0.0044604729969443596; D=-0.058869397444226292; E=-0.04138459079836742},
    @{Row=54; C=0.0060364711205798333; D=0.090310281688657928; E=0.11397291353038275},
    @{Row=55; C=0.010925590611306902; D=-0.036992122897494324; E=0.0058355860610350124}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
